$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-24 Friday" "2025-10-25 Saturday"

Replace-Text "87×20=" "26×25="
Replace-Text "54×76=" "92×98="
Replace-Text "64×69=" "79×36="
Replace-Text "82×60=" "16×76="
Replace-Text "73×16=" "50×93="
Replace-Text "85×70=" "39×66="
Replace-Text "88×92=" "69×88="
Replace-Text "38×59=" "14×25="
Replace-Text "29×12=" "56×79="
Replace-Text "82×15=" "53×13="
Replace-Text "52×44=" "14×86="
Replace-Text "27×29=" "92×75="
Replace-Text "24×63=" "66×23="
Replace-Text "86×33=" "80×42="
Replace-Text "83×88=" "86×59="
Replace-Text "56×35=" "22×76="
Replace-Text "28×86=" "97×44="
Replace-Text "87×93=" "63×54="
Replace-Text "30×93=" "64×33="
Replace-Text "46×53=" "35×44="
Replace-Text "33×81=" "37×25="
Replace-Text "76×99=" "17×71="
Replace-Text "66×79=" "76×37="
Replace-Text "47×16=" "11×56="
Replace-Text "19×56=" "63×99="
